$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, matching style of existing headers (B1:H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-12
$data = @(
    @(8, 9),
    @(5, 7),
    @(6, 7),
    @(7, 9),
    @(2, 4),
    @(1, 2),
    @(5, 7),
    @(1, 4),
    @(5, 7),
    @(6, 7),
    @(5, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
